$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'8.30%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'50.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'20.45%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.378"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'7.52%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08103"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'7.73%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.576"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.42%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.640"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.25%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.095"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'18.24%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1317"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'10.32%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1954"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'6.80%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09507"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.50%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04547"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11.20%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D14").Value = "'0.001324"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.79%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005792"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-3.68%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.382"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.62%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.431"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.23%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'2.34%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'8.170"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.12%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1398"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.42%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2927"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-11.33%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04303"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.79%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D24").Value = "'0.004268"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'9.60%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'9.58%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003717"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.17%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E38").Value = "'14.97%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05536"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'6.44%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006290"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.22%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007773"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.41%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1444"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'9.09%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007684"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.77%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008806"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'19.20%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3526"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'19.48%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006814"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'5.71%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.06059"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'76.92%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003995"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.93%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("E51").Style = "Normal"
